# Apply the "invert binary tree and dynamic programming" changes to the
# AlgoExpert tracking sheet:
#   - Row 21 (BST Validation): fill in empty Pain Points/Date Revisited/
#     Date Attempted First/Insight cells (still blank, but styled like the
#     rest of the row) so the row extends to column I.
#   - Row 22 (PreOrder, InOrder and PostOrder): add the Pain Points note
#     about pre/in/post-order traversal, plus the same styled blanks for
#     Date Revisited/Date Attempted First/Insight.
#   - Row 23 (new): "Invert Binary Tree" problem row with a Pain Points
#     note about traversing the tree.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style index 1 (pink-filled rows) is used throughout rows 17-22; grab its
# fill color from an existing cell in that block so the new cells match
# formatting.
$rowFill = $ws.Range("A21").Interior.Color

# --- Row 21: BST Validation -------------------------------------------
$ws.Range("F21:I21").Interior.Color = $rowFill

# --- Row 22: PreOrder, InOrder and PostOrder ---------------------------
$ws.Range("F22").Value = "Pre-order means append to array first, then traverse left, then right. Post-order append last, in Order means append in the middle"
$ws.Range("F22:I22").Interior.Color = $rowFill

# --- Row 23: Invert Binary Tree (new row) ------------------------------
$ws.Range("A23:F23").Interior.Color = $rowFill
$ws.Range("A23").Value = "Invert Binary Tree"
$ws.Range("B23").Value = "40 mins"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = "90 mins"
$ws.Range("E23").Value = "Failure"
$ws.Range("F23").Value = "How to traverse through the tree, "

$ws.Range("G23").Select()

$wb.Save()
